$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete shipment rows (old rows 2 and 3 - TCKU/MRSU), keeping
# only the MRKU row (old row 4), which shifts up to become the new row 2.
# This also updates the B2 value to the MRKU shipment's number (701925) and
# refreshes the used range / shared-string table automatically.
$ws.Rows("2:3").Delete()
